$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 21:18:40"
$ws.Range("N2").Value = "-2.5 °C 20:37 TU"
$ws.Range("E3").Value = "2026-02-20 21:18:42"
$ws.Range("O3").Value = "-5.2 °C"
$ws.Range("E4").Value = "2026-02-20 21:18:45"
$ws.Range("H4").Value = "'58%"
$ws.Range("J4").Value = "1022.7 hPa"
$ws.Range("N4").Value = "4.8 °C 20:58 TU"
$ws.Range("O4").Value = "10.1 °C"
$ws.Range("E5").Value = "2026-02-20 21:18:47"
$ws.Range("E6").Value = "2026-02-20 21:18:50"
$ws.Range("H6").Value = "'68%"
$ws.Range("J6").Value = "1022.7 hPa"
$ws.Range("E7").Value = "2026-02-20 21:18:53"
$ws.Range("H7").Value = "'48%"
$ws.Range("J7").Value = "1022.6 hPa"
$ws.Range("E8").Value = "2026-02-20 21:18:55"
$ws.Range("J8").Value = "1022.9 hPa"
$ws.Range("E9").Value = "2026-02-20 21:18:58"
$ws.Range("O9").Value = "13.3 °C"
$ws.Range("E10").Value = "2026-02-20 21:19:00"
$ws.Range("O10").Value = "7.6 °C"
$ws.Range("E11").Value = "2026-02-20 21:19:03"
$ws.Range("H11").Value = "'33%"
$ws.Range("O11").Value = "9.3 °C"
$ws.Range("E12").Value = "2026-02-20 21:19:06"
$ws.Range("H12").Value = "'48%"
$ws.Range("E13").Value = "2026-02-20 21:19:08"
$ws.Range("J13").Value = "1023.8 hPa"
$ws.Range("N13").Value = "0.0 °C 20:41 TU"
$ws.Range("O13").Value = "6.3 °C"
$ws.Range("E14").Value = "2026-02-20 21:19:11"
$ws.Range("H14").Value = "'58%"
$ws.Range("N14").Value = "6.3 °C 20:31 TU"
$ws.Range("O14").Value = "11.9 °C"
$ws.Range("E15").Value = "2026-02-20 21:19:13"
$ws.Range("E16").Value = "2026-02-20 21:19:16"
$ws.Range("O16").Value = "-3.1 °C"
$ws.Range("E17").Value = "2026-02-20 21:19:18"
$ws.Range("O17").Value = "3.1 °C"
$ws.Range("E18").Value = "2026-02-20 21:19:21"
$ws.Range("J18").Value = "1023.0 hPa"
$ws.Range("O18").Value = "7.9 °C"
$ws.Range("E19").Value = "2026-02-20 21:19:24"
$ws.Range("E20").Value = "2026-02-20 21:19:25"
$ws.Range("E21").Value = "2026-02-20 21:19:26"
$ws.Range("H21").Value = "'37%"
$ws.Range("J21").Value = "1022.8 hPa"
$ws.Range("O21").Value = "9.3 °C"
$ws.Range("E22").Value = "2026-02-20 21:19:27"
$ws.Range("O22").Value = "-4.0 °C"
$ws.Range("E23").Value = "2026-02-20 21:19:29"
$ws.Range("O23").Value = "-4.7 °C"
$ws.Range("E24").Value = "2026-02-20 21:19:31"
$ws.Range("E25").Value = "2026-02-20 21:19:34"
$ws.Range("E26").Value = "2026-02-20 21:19:36"
$ws.Range("J26").Value = "1021.9 hPa"
$ws.Range("E27").Value = "2026-02-20 21:19:39"
$ws.Range("O27").Value = "-0.7 °C"
$ws.Range("E28").Value = "2026-02-20 21:19:41"
$ws.Range("H28").Value = "'65%"
$ws.Range("J28").Value = "1023.1 hPa"
$ws.Range("O28").Value = "7.0 °C"
$ws.Range("E29").Value = "2026-02-20 21:19:44"
$ws.Range("H29").Value = "'69%"
$ws.Range("O29").Value = "9.9 °C"
$ws.Range("E30").Value = "2026-02-20 21:19:46"
$ws.Range("H30").Value = "'59%"
$ws.Range("J30").Value = "1022.4 hPa"
$ws.Range("O30").Value = "11.0 °C"
$ws.Range("E31").Value = "2026-02-20 21:19:49"
$ws.Range("J31").Value = "1021.7 hPa"
$ws.Range("O31").Value = "10.8 °C"
$ws.Range("E32").Value = "2026-02-20 21:19:52"
$ws.Range("H32").Value = "'83%"
$ws.Range("E33").Value = "2026-02-20 21:19:54"
$ws.Range("H33").Value = "'42%"
$ws.Range("J33").Value = "1023.2 hPa"
$ws.Range("N33").Value = "2.2 °C 20:51 TU"
$ws.Range("E34").Value = "2026-02-20 21:19:57"
$ws.Range("O34").Value = "1.0 °C"
$ws.Range("E35").Value = "2026-02-20 21:19:59"
$ws.Range("J35").Value = "1026.9 hPa"
$ws.Range("E36").Value = "2026-02-20 21:20:02"
$ws.Range("J36").Value = "1022.6 hPa"
$ws.Range("E37").Value = "2026-02-20 21:20:04"
$ws.Range("H37").Value = "'67%"
$ws.Range("J37").Value = "1024.6 hPa"
$ws.Range("O37").Value = "4.7 °C"
$ws.Range("E38").Value = "2026-02-20 21:20:07"
$ws.Range("K38").Value = "9.1 MJ/m2"
$ws.Range("O38").Value = "8.8 °C"
$ws.Range("E39").Value = "2026-02-20 21:20:09"
$ws.Range("O39").Value = "-2.5 °C"
$ws.Range("E40").Value = "2026-02-20 21:20:12"
$ws.Range("H40").Value = "'37%"
$ws.Range("J40").Value = "1023.6 hPa"
$ws.Range("O40").Value = "10.4 °C"
$ws.Range("E41").Value = "2026-02-20 21:20:14"
$ws.Range("E42").Value = "2026-02-20 21:20:16"
$ws.Range("O42").Value = "10.4 °C"
$ws.Range("E43").Value = "2026-02-20 21:20:19"
$ws.Range("E44").Value = "2026-02-20 21:20:21"
$ws.Range("K44").Value = "9.9 MJ/m2"
$ws.Range("O44").Value = "-4.6 °C"
$ws.Range("E45").Value = "2026-02-20 21:20:23"
$ws.Range("J45").Value = "1029.7 hPa"
$ws.Range("E46").Value = "2026-02-20 21:20:26"
$ws.Range("H46").Value = "'56%"
$ws.Range("J46").Value = "1026.5 hPa"
